# Adapt the two "old"/"new" header blocks to be named after the actual
# input file versions (FV2410 / FV2504), wrap the header+data range in a
# worksheet Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 21 header cells in row 1 (A1:U1).
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range (A1:U58) into an Excel table ("Table1") so the new
#    header names are also reflected in the table column definitions.
$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
